$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 152, shifting existing rows 152:275 down to 153:276
# (this is a new weekly price observation being added into the series).
$ws.Rows("152:152").Insert()

$ws.Range("A152").Value = 10
$ws.Range("B152").Value = "Vega Modelo de Temuco"
$ws.Range("C152").Value = "La Araucanía"
$ws.Range("D152").Value = 44658
$ws.Range("E152").Value = 9
$ws.Range("F152").Value = 100112017
$ws.Range("G152").Value = "Apio"
$ws.Range("H152").Value = "Americana (o)"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 215
$ws.Range("K152").Value = 9000
$ws.Range("L152").Value = 10000
$ws.Range("M152").Value = 9419
$ws.Range("N152").Value = "`$/docena de matas"
$ws.Range("O152").Value = "Provincia del Elquí"
$ws.Range("P152").Value = 1570
$ws.Range("Q152").Value = 6
$ws.Range("R152").Value = "Hortaliza"
